$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39 / Row 40 swap: VeChain <-> TrustWalletToken (Coin name + Link columns) ---
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'

# --- Price column (D): values are stored as text in the sheet, even when they look
# numeric (e.g. "1.002", "26.761.92"). A leading apostrophe forces Excel to keep the
# literal text instead of auto-converting it to a number, and resetting the Style
# afterwards clears the quote-prefix formatting flag Excel attaches to the cell. ---
$ws.Range("D2").Value = '''26.761.92'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = '''1.868.99'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Value = '''300.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = '''0.5323'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = '''0.3728'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = '''0.07146'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = '''21.45'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = '''0.8859'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = '''0.08154'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = '''1.879.30'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = '''92.38'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = '''5.293'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = '''1.002'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = '''14.84'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = '''0.000008491'
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Value = '''26.796.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = '''4.975'
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = '''6.378'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = '''2.285'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = '''145.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = '''1.731'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = '''18.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = '''113.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = '''4.695'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = '''4.626'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = '''0.09105'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = '''0.8104'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = '''0.05018'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = '''1.171'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = '''2.944'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = '''0.6120'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = '''2.648'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = '''3.182'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = '''0.01942'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = '''1.068'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = '''0.5301'
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Value = '''8.721'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = '''115.69'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = '''0.1491'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = '''1.000'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = '''1.649'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = '''9.950'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = '''37.30'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = '''0.06060'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = '''62.12'
$ws.Range("D51").Style = "Normal"

# --- Volume(1h) column (E): padded percentage strings (e.g. "  -1.86%  ") are never
# auto-converted to numbers by Excel, so a plain text assignment is sufficient. ---
$ws.Range("E2").Value = '  -1.86%  '
$ws.Range("E3").Value = '  -2.07%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E5").Value = '  -2.36%  '
$ws.Range("E7").Value = '  +1.11%  '
$ws.Range("E8").Value = '  -2.18%  '
$ws.Range("E9").Value = '  -2.02%  '
$ws.Range("E10").Value = '  -2.40%  '
$ws.Range("E11").Value = '  -1.83%  '
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("E13").Value = '  +29.55%  '
$ws.Range("E14").Value = '  -3.80%  '
$ws.Range("E15").Value = '  -1.30%  '
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("E18").Value = '  -2.03%  '
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("E21").Value = '  -2.60%  '
$ws.Range("E22").Value = '  -1.82%  '
$ws.Range("E23").Value = '  -2.18%  '
$ws.Range("E24").Value = '  -0.80%  '
$ws.Range("E25").Value = '  -2.84%  '
$ws.Range("E26").Value = '  -0.66%  '
$ws.Range("E27").Value = '  -1.32%  '
$ws.Range("E28").Value = '  -2.63%  '
$ws.Range("E29").Value = '  -3.26%  '
$ws.Range("E30").Value = '  -4.66%  '
$ws.Range("E31").Value = '  -1.65%  '
$ws.Range("E32").Value = '  -2.36%  '
$ws.Range("E34").Value = '  -4.74%  '
$ws.Range("E36").Value = '  +5.56%  '
$ws.Range("E37").Value = '  -2.52%  '
$ws.Range("E38").Value = '  -5.04%  '
$ws.Range("E39").Value = '  -2.92%  '
$ws.Range("E40").Value = '  -0.69%  '
$ws.Range("E41").Value = '  +7.57%  '
$ws.Range("E42").Value = '  -1.67%  '
$ws.Range("E43").Value = '  -5.51%  '
$ws.Range("E44").Value = '  -0.83%  '
$ws.Range("E45").Value = '  -2.25%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("E48").Value = '  -2.55%  '
$ws.Range("E49").Value = '  -4.16%  '
$ws.Range("E50").Value = '  -2.32%  '
$ws.Range("E51").Value = '  -3.73%  '
